$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 202.92857
$ws.Range("I33").Value = 246.54546
$ws.Range("J33").Value = 43
$ws.Range("K33").Value = 246.54546
$ws.Range("L33").Value = 43
$ws.Range("M33").Value = -17.54545999999999
$ws.Range("N33").Value = -501
$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872
$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360
$ws.Range("H80").Value = 9167.789000000001
$ws.Range("J80").Value = 13099.917
$ws.Range("L80").Value = 39299.751
$ws.Range("N80").Value = -41295.751
$ws.Range("H83").Value = 9167.789000000001
$ws.Range("J83").Value = 13099.917
$ws.Range("L83").Value = 117899.253
$ws.Range("N83").Value = -127883.253
$ws.Range("H110").Value = 52776.715
$ws.Range("J110").Value = 52776.715
$ws.Range("L110").Value = 52776.715
$ws.Range("N110").Value = -60956.715
$ws.Range("H131").Value = 453.25
$ws.Range("I131").Value = 104.333336
$ws.Range("J131").Value = 1500
$ws.Range("K131").Value = 313.000008
$ws.Range("L131").Value = 4500
$ws.Range("M131").Value = 4726.999992
$ws.Range("N131").Value = -14580
$ws.Range("H133").Value = 76650.60000000001
$ws.Range("J133").Value = 76650.60000000001
$ws.Range("L133").Value = 76650.60000000001
$ws.Range("N133").Value = -86770.60000000001
$ws.Range("H134").Value = 94216.664
$ws.Range("J134").Value = 94216.664
$ws.Range("L134").Value = 94216.664
$ws.Range("N134").Value = -104356.664
$ws.Range("H136").Value = 96495.836
$ws.Range("J136").Value = 96495.836
$ws.Range("L136").Value = 96495.836
$ws.Range("N136").Value = -106695.836
$ws.Range("H139").Value = 99406
$ws.Range("J139").Value = 99406
$ws.Range("L139").Value = 99406
$ws.Range("N139").Value = -109686

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 46664.668
$ws.Range("J107").Value = 46664.668
$ws.Range("L107").Value = 46664.668
$ws.Range("N107").Value = -54344.668
$ws.Range("H130").Value = 88877
$ws.Range("J130").Value = 88877
$ws.Range("L130").Value = 88877
$ws.Range("N130").Value = -98917
$ws.Range("H131").Value = 99997
$ws.Range("J131").Value = 99997
$ws.Range("L131").Value = 99997
$ws.Range("N131").Value = -110077

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 99990
$ws.Range("J52").Value = 99990
$ws.Range("L52").Value = 99990
$ws.Range("N52").Value = -100516
$ws.Range("H54").Value = 3064.6
$ws.Range("I54").Value = 1738.4445
$ws.Range("K54").Value = 1738.4445
$ws.Range("M54").Value = -1254.4445
$ws.Range("H55").Value = 37121.5
$ws.Range("J55").Value = 37121.5
$ws.Range("L55").Value = 37121.5
$ws.Range("N55").Value = -37667.5
$ws.Range("H109").Value = 72996
$ws.Range("J109").Value = 72996
$ws.Range("L109").Value = 72996
$ws.Range("N109").Value = -75770
$ws.Range("H115").Value = 91624.375
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("H119").Value = 99990
$ws.Range("J119").Value = 99990
$ws.Range("L119").Value = 99990
$ws.Range("N119").Value = -109666
$ws.Range("H121").Value = 99990
$ws.Range("J121").Value = 99990
$ws.Range("L121").Value = 99990
$ws.Range("N121").Value = -103484
$ws.Range("H132").Value = 46554.71
$ws.Range("J132").Value = 46554.71
$ws.Range("L132").Value = 46554.71
$ws.Range("N132").Value = -56674.71
$ws.Range("H135").Value = 105567.71
$ws.Range("J135").Value = 105567.71
$ws.Range("L135").Value = 105567.71
$ws.Range("N135").Value = -115707.71
$ws.Range("H138").Value = 99999
$ws.Range("J138").Value = 99999
$ws.Range("L138").Value = 99999
$ws.Range("N138").Value = -110279
$ws.Range("H140").Value = 43433.934
$ws.Range("J140").Value = 43499.406
$ws.Range("L140").Value = 43499.406
$ws.Range("N140").Value = -53859.406

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2011.7646
$ws.Range("I31").Value = 1591.9354
$ws.Range("J31").Value = 6350
$ws.Range("K31").Value = 1591.9354
$ws.Range("L31").Value = 6350
$ws.Range("M31").Value = -1296.9354
$ws.Range("N31").Value = -6940
$ws.Range("H34").Value = 2011.7646
$ws.Range("I34").Value = 1591.9354
$ws.Range("J34").Value = 6350
$ws.Range("K34").Value = 1591.9354
$ws.Range("L34").Value = 6350
$ws.Range("M34").Value = -1389.9354
$ws.Range("N34").Value = -6754
$ws.Range("H114").Value = 39267.285
$ws.Range("J114").Value = 39267.285
$ws.Range("L114").Value = 39267.285
$ws.Range("N114").Value = -47945.285
$ws.Range("H118").Value = 64797.145
$ws.Range("J118").Value = 64797.145
$ws.Range("L118").Value = 64797.145
$ws.Range("N118").Value = -68111.14499999999
$ws.Range("H138").Value = 102243.43
$ws.Range("J138").Value = 108332.5
$ws.Range("L138").Value = 108332.5
$ws.Range("N138").Value = -118612.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 5187.625
$ws.Range("I7").Value = 70.2
$ws.Range("K7").Value = 210.6
$ws.Range("M7").Value = -98.60000000000002
$ws.Range("H92").Value = 315.83334
$ws.Range("J92").Value = 365
$ws.Range("L92").Value = 1095
$ws.Range("N92").Value = -3591
$ws.Range("H114").Value = 5928.467
$ws.Range("I114").Value = 728
$ws.Range("J114").Value = 7819.5454
$ws.Range("K114").Value = 2184
$ws.Range("L114").Value = 23458.6362
$ws.Range("M114").Value = 1070
$ws.Range("N114").Value = -29966.6362
$ws.Range("H121").Value = 2844
$ws.Range("J121").Value = 3812.25
$ws.Range("L121").Value = 11436.75
$ws.Range("N121").Value = -14056.75
$ws.Range("H137").Value = 6379.273
$ws.Range("I137").Value = 3892.375
$ws.Range("J137").Value = 13011
$ws.Range("K137").Value = 11677.125
$ws.Range("L137").Value = 39033
$ws.Range("M137").Value = -6577.125
$ws.Range("N137").Value = -49233

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 20666.9
$ws.Range("J93").Value = 20666.9
$ws.Range("L93").Value = 20666.9
$ws.Range("N93").Value = -24410.9
$ws.Range("H110").Value = 74221
$ws.Range("J110").Value = 74221
$ws.Range("L110").Value = 74221
$ws.Range("N110").Value = -82401
$ws.Range("H116").Value = 55872.5
$ws.Range("I116").Value = 33000
$ws.Range("K116").Value = 33000
$ws.Range("M116").Value = -28411
$ws.Range("H140").Value = 97330.11
$ws.Range("J140").Value = 98183.875
$ws.Range("L140").Value = 98183.875
$ws.Range("N140").Value = -108543.875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4184.4614
$ws.Range("I7").Value = 3318.5
$ws.Range("K7").Value = 3318.5
$ws.Range("M7").Value = -3206.5
$ws.Range("H118").Value = 50863.637
$ws.Range("J118").Value = 51950
$ws.Range("L118").Value = 51950
$ws.Range("N118").Value = -55264
$ws.Range("H122").Value = 66671104
$ws.Range("I122").Value = 71432984
$ws.Range("J122").Value = 50004516
$ws.Range("K122").Value = 214298952
$ws.Range("L122").Value = 150013548
$ws.Range("M122").Value = -214296502
$ws.Range("N122").Value = -150018448
$ws.Range("H123").Value = 79123
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("H126").Value = 4184.4614
$ws.Range("I126").Value = 3318.5
$ws.Range("K126").Value = 9955.5
$ws.Range("M126").Value = -7485.5
$ws.Range("H129").Value = 107231.5
$ws.Range("I129").Value = 76390
$ws.Range("K129").Value = 76390
$ws.Range("M129").Value = -71390
$ws.Range("H132").Value = 3922.8333
$ws.Range("I132").Value = 3562.375
$ws.Range("J132").Value = 4643.75
$ws.Range("K132").Value = 10687.125
$ws.Range("L132").Value = 13931.25
$ws.Range("M132").Value = -8157.125
$ws.Range("N132").Value = -18991.25
